$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: "Baseplate" renamed to "Tray 0" (quantity etc. unchanged)
$ws.Range("A11").Value = "Tray 0"

# Row 12 (Receiver Housing): quantity of 2 moves from Tray 1 (col B) to Tray 5 (col F)
$ws.Range("B12").Value = $null
$ws.Range("F12").Value = 2

# Row 13 (Transmitter Housing): quantity of 2 moves from Tray 5 (col F) to Tray 1 (col B)
$ws.Range("F13").Value = $null
$ws.Range("B13").Value = 2

# New rows: Made by / Edited by info
$ws.Range("A17").Value = "Made by"
$ws.Range("B17").Value = "S. Bartlett"
$ws.Range("A18").Value = "Edited by"
$ws.Range("B18").Value = "T. Nguyen"

$ws.Range("K15").Select()
